$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row 11 (CENTRO DE SALUD SAN CAMILO) - shifts old rows 11-43 down to 12-44
$ws.Rows.Item(11).Insert()

# Insert new row 30 (SERVICIOS DE SALUD IPS SURAMERICANA SA) - shifts rows 30-44 down to 31-45
$ws.Rows.Item(30).Insert()

# Update dimension-affecting header + data cells

$ws.Range("I1").Value = '6'

$ws.Range("I2").Value = 98

$ws.Range("I3").Value = 83

$ws.Range("I4").Value = 10

$ws.Range("I5").Value = 43

$ws.Range("I6").Value = 33

$ws.Range("I7").Value = 27

$ws.Range("I8").Value = 2

$ws.Range("I9").Value = 2

$ws.Range("I10").Value = 2

$ws.Range("A11").Value = '6600100332'
$ws.Range("B11").Value = '07'
$ws.Range("C11").Value = 'CENTRO DE SALUD SAN CAMILO'
$ws.Range("I11").Value = 1

$ws.Range("I12").Value = 3

$ws.Range("I13").Value = 2

$ws.Range("I14").Value = 1

$ws.Range("I18").Value = 1

$ws.Range("I19").Value = 2

$ws.Range("I20").Value = 26

$ws.Range("I22").Value = 166

$ws.Range("I23").Value = 7

$ws.Range("I24").Value = 4

$ws.Range("I25").Value = 6

$ws.Range("I26").Value = 17

$ws.Range("I27").Value = 4

$ws.Range("I28").Value = 13

$ws.Range("I29").Value = 55

$ws.Range("A30").Value = '6600102411'
$ws.Range("B30").Value = '02'
$ws.Range("C30").Value = 'SERVICIOS DE SALUD IPS SURAMERICANA SA'
$ws.Range("I30").Value = 1

$ws.Range("I31").Value = 62

$ws.Range("I32").Value = 24

$ws.Range("I33").Value = 151

$ws.Range("I34").Value = 89

$ws.Range("I35").Value = 162

$ws.Range("I36").Value = 6

$ws.Range("I37").Value = 106

$ws.Range("I38").Value = 1

$ws.Range("I39").Value = 0

$ws.Range("I40").Value = 3

$ws.Range("I42").Value = 39

$ws.Range("I43").Value = 3

$ws.Range("I44").Value = 12

$ws.Range("I45").Value = 42

